$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5: J5 takes over the old "kommentarstandarden" text (previously in H6),
# but is now explicitly colored black, which creates a distinct font/style.
$ws.Range("J5").Value = "Dokumentere funktioner i kode og med kommentarstandarden"
$ws.Range("J5").Font.Color = 0

# Row 6: I6's text is expanded to mention whitebox/blackbox testing.
$ws.Range("I6").Value = "Gennemføre en whitebox og blackbox afprøvning"

# Row 6: H6 gets brand-new content about signature files.
$ws.Range("H6").Value = "Kunne skrive og bruge en signaturfil"

# Move the active selection from H6 to G6.
[void]$ws.Range("G6").Select()
